# Update column G ("K" = strikeouts, renamed from Strike#) values for several rows
# per regen of save_data using K instead of Strike#.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 1
$ws.Range("G7").Value = 2
$ws.Range("G8").Value = 1
